$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): columns C/D change from the separate "Ouvrage (/)" /
# "Prestation (/)" pair to a combined "Ouvrage/Prestation" text column plus a
# new numeric "NbPrixRef" column (part of the new ImportBasePrixRef() service).
$ws.Range("C1").Value = "Ouvrage/Prestation"
$ws.Range("D1").Value = "NbPrixRef"

# Row 2: the old C2/D2 pair ("04_AAA_01_01_01" / "/") becomes a single
# combined reference in C2 plus a numeric ref count in D2.
$ws.Range("C2").Value = "04_AAA_01_01_01"
$ws.Range("D2").Value = 3

# Row 3: same restructuring.
$ws.Range("C3").Value = "04_AAA_01_01_01_01"
$ws.Range("D3").Value = 2

# Row 4: same restructuring.
$ws.Range("C4").Value = "04_AAA_01_01_01_02"
$ws.Range("D4").Value = 4

# Column widths followed the new (longer) text in C and the new (shorter)
# numeric content in D. ColumnWidth is expressed in characters; the stored
# OOXML width is ColumnWidth + 5/6, so back that off to land on the target
# stored widths of 19.83203125 / 9.5.
$ws.Columns("C").ColumnWidth = 18.998697916666668
$ws.Columns("D").ColumnWidth = 8.666666666666666

# Selection moved to E4.
$ws.Range("E4").Select()
